# Update the cryptos price/volume table (columns D and E, rows 2-51) with
# freshly scraped values. Column D/E cells are stored as text in the
# workbook (prices like "42.040.84" use '.' as a thousands separator, and
# the Volume column is a padded "  +x.xx%  " string), so a leading "'" is
# used on any new value that Excel's auto-detection would otherwise treat
# as a plain number (which would silently drop significant trailing
# zeros, e.g. "3.60" -> 3.6, or convert the cell to a numeric type).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.981.64"
$ws.Range("E2").Value = "  -0.57%  "

$ws.Range("D3").Value = "2.210.15"
$ws.Range("E3").Value = "  -1.45%  "

$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("D5").Value = "'240.22"
$ws.Range("E5").Value = "  -2.58%  "

$ws.Range("D6").Value = "'0.623"
$ws.Range("E6").Value = "  -0.83%  "

$ws.Range("D7").Value = "'72.76"
$ws.Range("E7").Value = "  -1.91%  "

$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("D9").Value = "'0.600"
$ws.Range("E9").Value = "  -2.77%  "

$ws.Range("D10").Value = "'41.97"
$ws.Range("E10").Value = "  -1.00%  "

$ws.Range("D11").Value = "'0.0946"
$ws.Range("E11").Value = "  +0.97%  "

$ws.Range("D12").Value = "'7.01"
$ws.Range("E12").Value = "  -1.11%  "

$ws.Range("E13").Value = "  -0.13%  "

$ws.Range("D14").Value = "2.541.73"
$ws.Range("E14").Value = "  -1.38%  "

$ws.Range("D15").Value = "'14.11"
$ws.Range("E15").Value = "  -2.67%  "

$ws.Range("D16").Value = "'0.831"
$ws.Range("E16").Value = "  -2.36%  "

$ws.Range("D17").Value = "2.223.52"
$ws.Range("E17").Value = "  -0.73%  "

$ws.Range("D18").Value = "41.855.54"
$ws.Range("E18").Value = "  -0.57%  "

$ws.Range("E19").Value = "  +8.89%  "

$ws.Range("D20").Value = "'72.82"
$ws.Range("E20").Value = "  +1.01%  "

$ws.Range("D21").Value = "'6.09"
$ws.Range("E21").Value = "  -0.98%  "

$ws.Range("D22").Value = "'10.13"
$ws.Range("E22").Value = "  +14.24%  "

$ws.Range("D23").Value = "'228.61"
$ws.Range("E23").Value = "  -1.26%  "

$ws.Range("E24").Value = "  -7.06%  "

$ws.Range("D25").Value = "'11.57"
$ws.Range("E25").Value = "  +2.45%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("D27").Value = "'3.60"
$ws.Range("E27").Value = "  -0.70%  "

$ws.Range("D28").Value = "'2.25"
$ws.Range("E28").Value = "  -1.66%  "

$ws.Range("E29").Value = "  +1.10%  "

$ws.Range("D30").Value = "'166.76"
$ws.Range("E30").Value = "  -1.32%  "

$ws.Range("D31").Value = "'20.50"
$ws.Range("E31").Value = "  -0.68%  "

$ws.Range("D32").Value = "'5.62"
$ws.Range("E32").Value = "  +7.52%  "

$ws.Range("D33").Value = "'0.0784"
$ws.Range("E33").Value = "  -4.23%  "

$ws.Range("E34").Value = "  -0.42%  "

$ws.Range("D35").Value = "'28.54"
$ws.Range("E35").Value = "  -6.36%  "

$ws.Range("D36").Value = "'0.109"
$ws.Range("E36").Value = "  -8.37%  "

$ws.Range("D37").Value = "'4.21"
$ws.Range("E37").Value = "  -5.36%  "

$ws.Range("D38").Value = "'0.0297"
$ws.Range("E38").Value = "  -4.45%  "

$ws.Range("D39").Value = "'13.22"
$ws.Range("E39").Value = "  -2.41%  "

$ws.Range("D40").Value = "'64.97"
$ws.Range("E40").Value = "  +4.70%  "

$ws.Range("D41").Value = "'2.10"
$ws.Range("E41").Value = "  -3.85%  "

$ws.Range("D42").Value = "'5.59"
$ws.Range("E42").Value = "  -2.93%  "

$ws.Range("D43").Value = "'0.196"
$ws.Range("E43").Value = "  -4.07%  "

$ws.Range("D44").Value = "'8.64"
$ws.Range("E44").Value = "  -0.61%  "

$ws.Range("D45").Value = "'103.36"
$ws.Range("E45").Value = "  -3.04%  "

$ws.Range("E46").Value = "  -2.24%  "

$ws.Range("D47").Value = "'2.35"
$ws.Range("E47").Value = "  +3.53%  "

$ws.Range("D48").Value = "'1.10"
$ws.Range("E48").Value = "  -1.15%  "

$ws.Range("E49").Value = "  -0.18%  "

$ws.Range("E50").Value = "  -0.11%  "

$ws.Range("D51").Value = "2.418.57"
$ws.Range("E51").Value = "  -2.23%  "
